# Applies corrected IFRS financial figures to rows 2-9 of the company_list sheet.
# (commit: "error solve ifrs list")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value corrections (cell address -> new value) ---
$updates = @(
    @("D2", 2786),
    @("E2", -404),
    @("F2", -404),
    @("G2", -1060),
    @("H2", -1060),
    @("I2", -1060),
    @("K2", 5256),
    @("L2", 5915),
    @("M2", -659),
    @("N2", -661),
    @("O2", 2),
    @("P2", 414),
    @("Q2", 19),
    @("R2", 49),
    @("S2", -9),
    @("T2", 0),
    @("U2", 18),
    @("V2", 92),
    @("W2", -14.51),
    @("X2", -38.04),
    @("Y2", 555.95),
    @("Z2", -21.55),
    @("AA2", -897.26),
    @("AB2", -256.79),
    @("AC2", -531767),
    @("AD2", -0.22),
    @("AE2", -321659),
    @("AF2", -0.36),
    @("AG2", 0),
    @("AH2", 0),
    @("AI2", 0),
    @("AJ2", 205448),
    @("D3", 3024),
    @("E3", -509),
    @("F3", -509),
    @("G3", -770),
    @("H3", -770),
    @("I3", -770),
    @("K3", 2392),
    @("L3", 1998),
    @("M3", 394),
    @("N3", 394),
    @("P3", 490),
    @("Q3", 73),
    @("R3", -278),
    @("S3", 268),
    @("T3", 0),
    @("U3", 72),
    @("V3", 92),
    @("W3", -16.84),
    @("X3", -25.47),
    @("Y3", 578.25),
    @("Z3", -20.14),
    @("AA3", 506.58),
    @("AB3", 60.02),
    @("AC3", -49479),
    @("AD3", -0.44),
    @("AE3", 4026),
    @("AF3", 5.38),
    @("AG3", 0),
    @("AH3", 0),
    @("AI3", 0),
    @("AJ3", 9796455),
    @("D4", 1836),
    @("E4", -93),
    @("F4", -93),
    @("G4", -25),
    @("H4", -21),
    @("I4", -21),
    @("K4", 1759),
    @("L4", 1352),
    @("M4", 407),
    @("N4", 407),
    @("P4", 490),
    @("Q4", -46),
    @("R4", 334),
    @("S4", -293),
    @("T4", 3),
    @("U4", -49),
    @("V4", 92),
    @("W4", -5.07),
    @("X4", -1.17),
    @("Y4", -5.37),
    @("Z4", -1.03),
    @("AA4", 332.68),
    @("AB4", 56.93),
    @("AC4", -219),
    @("AD4", -33.18),
    @("AE4", 4153),
    @("AF4", 1.75),
    @("AG4", 0),
    @("AH4", 0),
    @("AI4", 0),
    @("AJ4", 9796455),
    @("D5", 1156),
    @("E5", 18),
    @("F5", 18),
    @("G5", 16),
    @("H5", 17),
    @("I5", 17),
    @("K5", 1664),
    @("L5", 1231),
    @("M5", 433),
    @("N5", 433),
    @("P5", 491),
    @("Q5", 13),
    @("R5", 21),
    @("S5", 0),
    @("T5", 12),
    @("U5", 1),
    @("V5", 92),
    @("W5", 1.58),
    @("X5", 1.47),
    @("Y5", 4.05),
    @("Z5", 0.99),
    @("AA5", 284.24),
    @("AB5", 60.8),
    @("AC5", 173),
    @("AD5", 34.76),
    @("AE5", 4414),
    @("AF5", 1.37),
    @("AG5", 0),
    @("AH5", 0),
    @("AI5", 0),
    @("AJ5", 9819185),
    @("D6", 1572),
    @("E6", 111),
    @("F6", 111),
    @("G6", 57),
    @("H6", 57),
    @("I6", 57),
    @("K6", 1489),
    @("L6", 998),
    @("M6", 491),
    @("N6", 491),
    @("P6", 491),
    @("Q6", 50),
    @("R6", 24),
    @("S6", 0),
    @("T6", 3),
    @("U6", 47),
    @("V6", 92),
    @("W6", 7.09),
    @("X6", 3.64),
    @("Y6", 12.37),
    @("Z6", 3.62),
    @("AA6", 203.19),
    @("AB6", 75.62),
    @("AC6", 582),
    @("AD6", 20.46),
    @("AE6", 5002),
    @("AF6", 2.38),
    @("AG6", 0),
    @("AH6", 0),
    @("AI6", 0),
    @("AJ6", 9827396)
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# --- Cells that no longer hold data and must be cleared ---
$clears = @(
    "J2", "J3", "O3", "J4", "O4", "J5", "O5", "D7",
    "E7", "G7", "H7", "I7", "K7", "L7", "M7", "N7",
    "P7", "Q7", "R7", "S7", "T7", "U7", "W7", "X7",
    "Y7", "Z7", "AA7", "AC7", "AD7", "AE7", "AF7", "AG7",
    "AH7", "AI7", "D8", "E8", "G8", "H8", "I8", "K8",
    "L8", "M8", "N8", "P8", "Q8", "R8", "S8", "T8",
    "U8", "W8", "X8", "Y8", "Z8", "AA8", "AC8", "AD8",
    "AE8", "AF8", "AG8", "AH8", "AI8", "D9", "E9", "G9",
    "H9", "I9", "K9", "L9", "M9", "N9", "P9", "Q9",
    "R9", "S9", "T9", "U9", "W9", "X9", "Y9", "Z9",
    "AA9", "AC9", "AD9", "AE9", "AF9", "AG9", "AH9", "AI9"
)

foreach ($addr in $clears) {
    $ws.Range($addr).ClearContents()
}

Write-Output "Applied ifrs list corrections ($($updates.Count) updated, $($clears.Count) cleared)"
